# New translations s.xlsx (Multilingual)
# Row 3 of the "Worksheet" sheet previously held placeholder/test data
# ("6464" as the identifier, "Pug" as the English label). Replace it with
# the real identifier and the correct "Save as..." translation entry
# (the other languages in that row were already correct).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Worksheet")

$ws.Range("A3").Value = 533
$ws.Range("B3").Value = "Save as..."
